# Generate Report for Handoff
#
# A new handoff round ran for the four files that were previously queued
# with "low" priority (6b6d8154, a745513b, c8113dad, ebe09e16). Their
# priority is now "ht" and their "Latest Handoff Datetime" is refreshed,
# for both locale sheets (zh-cn, de-de), rows 4-7.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Priority column (E): "low" -> "ht" for rows 4-7 in both locale sheets.
foreach ($row in 4..7) {
    $zhcn.Range("E$row").Value = "ht"
    $dede.Range("E$row").Value = "ht"
}

# Latest Handoff Datetime column (H): refresh timestamp for the same rows.
foreach ($row in 4..7) {
    $zhcn.Range("H$row").Value = "2016-08-21 22:42:14"
    $dede.Range("H$row").Value = "2016-08-21 22:42:18"
}

# Overview sheet tracks the latest (de-de) handoff generation date per file.
foreach ($row in 4..7) {
    $overview.Range("G$row").Value = "2016-08-21 22:42:18"
}
